# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Thu Jun 20 16:33:06 UTC 2024 with GitHub Actions"
#
# Updates Price (D) / Volume(1h) (E) values for most rows, and for a few
# rows the ranking reshuffled so the Coin name (B) + Link (C) also swap
# with the neighboring row (28/29, 34/35, 43/44).
#
# Price/Volume columns are stored as literal text in the sheet (e.g. "1.00",
# "  -0.61%  "), so NumberFormat is forced to "@" (Text) before writing any
# new Price value that would otherwise be auto-parsed as a number by Excel,
# keeping the cell a text cell exactly like the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.653.69'
$ws.Range('E2').Value = '  -0.61%  '

# Row 3
$ws.Range('D3').Value = '3.504.59'
$ws.Range('E3').Value = '  -1.15%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.76'
$ws.Range('E5').Value = '  -0.20%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.27'
$ws.Range('E6').Value = '  -3.62%  '

# Row 7
$ws.Range('D7').Value = '3.502.12'
$ws.Range('E7').Value = '  -1.19%  '

# Row 8
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -0.50%  '

# Row 10
$ws.Range('E10').Value = '  -0.29%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.10'
$ws.Range('E11').Value = '  +2.09%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.380'
$ws.Range('E12').Value = '  -1.05%  '

# Row 13
$ws.Range('D13').Value = '4.095.68'
$ws.Range('E13').Value = '  -1.26%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  -0.09%  '

# Row 15
$ws.Range('E15').Value = '  -1.65%  '

# Row 16
$ws.Range('E16').Value = '  +0.04%  '

# Row 17
$ws.Range('D17').Value = '3.503.26'
$ws.Range('E17').Value = '  -1.22%  '

# Row 18
$ws.Range('D18').Value = '64.689.44'
$ws.Range('E18').Value = '  -0.27%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.06'
$ws.Range('E19').Value = '  -1.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.28'
$ws.Range('E20').Value = '  +0.16%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.64'
$ws.Range('E21').Value = '  -2.98%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.83'
$ws.Range('E22').Value = '  +0.25%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.571'
$ws.Range('E23').Value = '  -0.62%  '

# Row 24
$ws.Range('D24').Value = '3.643.49'
$ws.Range('E24').Value = '  -1.01%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.66'
$ws.Range('E25').Value = '  -0.43%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.03%  '

# Row 27
$ws.Range('E27').Value = '  -3.32%  '

# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.59'
$ws.Range('E28').Value = '  -1.61%  '

# Row 29
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.57'
$ws.Range('E29').Value = '  +14.75%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.09%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.27'
$ws.Range('E31').Value = '  -0.24%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.25'
$ws.Range('E32').Value = '  +0.44%  '

# Row 33
$ws.Range('D33').Value = '3.508.34'
$ws.Range('E33').Value = '  -1.34%  '

# Row 34
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.95'
$ws.Range('E34').Value = '  +0.70%  '

# Row 35
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.04%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.143'
$ws.Range('E36').Value = '  -0.80%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '168.05'
$ws.Range('E37').Value = '  -0.78%  '

# Row 38
$ws.Range('E38').Value = '  -0.63%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.03'
$ws.Range('E39').Value = '  +1.61%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.77'
$ws.Range('E40').Value = '  -1.44%  '

# Row 41
$ws.Range('E41').Value = '  +0.64%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.815'
$ws.Range('E42').Value = '  -0.64%  '

# Row 43
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.53'
$ws.Range('E43').Value = '  -0.50%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.10%  '

# Row 45
$ws.Range('E45').Value = '  +1.43%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.16'
$ws.Range('E46').Value = '  -5.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.37'
$ws.Range('E47').Value = '  -1.22%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.63'
$ws.Range('E48').Value = '  -2.27%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.86'
$ws.Range('E49').Value = '  -0.23%  '

# Row 50
$ws.Range('D50').Value = '2.375.03'
$ws.Range('E50').Value = '  -2.29%  '

# Row 51
$ws.Range('E51').Value = '  +0.65%  '
